$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Semestre ideal:" value (row 9, columns B and C) from "EP-3" to "EF-7,EP-3"
$ws.Range("B9").Value = "EF-7,EP-3"
$ws.Range("C9").Value = "EF-7,EP-3"

# Remove the trailing "Requisitos:" row (22) and its value row (23)
$ws.Rows("22:23").Delete()
